$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("H62").Value = 25002908
$ws.Range("I62").Value = 50002884
$ws.Range("J62").Value = 2932.3
$ws.Range("K62").Value = 50002884
$ws.Range("L62").Value = 2932.3
$ws.Range("M62").Value = -50002260
$ws.Range("N62").Value = -4180.3
$ws.Range("H65").Value = 25002908
$ws.Range("I65").Value = 50002884
$ws.Range("J65").Value = 2932.3
$ws.Range("K65").Value = 250014420
$ws.Range("L65").Value = 14661.5
$ws.Range("M65").Value = -250011300
$ws.Range("N65").Value = -20901.5
$ws.Range("H100").Value = 1155.4615
$ws.Range("I100").Value = 1039.375
$ws.Range("J100").Value = 1341.2
$ws.Range("K100").Value = 1039.375
$ws.Range("L100").Value = 1341.2
$ws.Range("M100").Value = -498.375
$ws.Range("N100").Value = -2423.2
$ws.Range("H132").Value = 1723.7872
$ws.Range("I132").Value = 1720.4
$ws.Range("K132").Value = 5161.200000000001
$ws.Range("M132").Value = -2631.200000000001
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3225
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 3225
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 277.6
$ws.Range("I22").Value = 277.6
$ws.Range("K22").Value = 277.6
$ws.Range("M22").Value = 72.39999999999998
$ws.Range("H134").Value = 3245.6956
$ws.Range("I134").Value = 989.1539
$ws.Range("J134").Value = 6179.2
$ws.Range("K134").Value = 2967.4617
$ws.Range("L134").Value = 18537.6
$ws.Range("M134").Value = -432.4616999999998
$ws.Range("N134").Value = -23607.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2825.25
$ws.Range("I17").Value = 499.5
$ws.Range("J17").Value = 5151
$ws.Range("K17").Value = 1498.5
$ws.Range("L17").Value = 15453
$ws.Range("M17").Value = -1329.5
$ws.Range("N17").Value = -15791
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H50").Value = 66813.336
$ws.Range("J50").Value = 111312.22
$ws.Range("L50").Value = 333936.66
$ws.Range("N50").Value = -334898.66
$ws.Range("H53").Value = 66813.336
$ws.Range("J53").Value = 111312.22
$ws.Range("L53").Value = 333936.66
$ws.Range("N53").Value = -334898.66
$ws.Range("H59").Value = 1359.5834
$ws.Range("I59").Value = 278.75
$ws.Range("J59").Value = 1900
$ws.Range("K59").Value = 836.25
$ws.Range("L59").Value = 5700
$ws.Range("M59").Value = -296.25
$ws.Range("N59").Value = -6780
$ws.Range("H68").Value = 913.871
$ws.Range("I68").Value = 790.4
$ws.Range("J68").Value = 1073.9259
$ws.Range("K68").Value = 2371.2
$ws.Range("L68").Value = 3221.7777
$ws.Range("M68").Value = -1560.2
$ws.Range("N68").Value = -4843.7777
$ws.Range("H71").Value = 913.871
$ws.Range("I71").Value = 790.4
$ws.Range("J71").Value = 1073.9259
$ws.Range("K71").Value = 7113.599999999999
$ws.Range("L71").Value = 9665.3331
$ws.Range("M71").Value = -3057.599999999999
$ws.Range("N71").Value = -17777.3331
$ws.Range("H131").Value = 2427.0588
$ws.Range("I131").Value = 694.1429000000001
$ws.Range("J131").Value = 2625.918
$ws.Range("K131").Value = 2082.4287
$ws.Range("L131").Value = 7877.754000000001
$ws.Range("M131").Value = 2957.5713
$ws.Range("N131").Value = -17957.754

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 292.4
$ws.Range("I107").Value = 292.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 292.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1627.6
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 67334390
$ws.Range("I22").Value = 77693340
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 77693340
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -77693045
$ws.Range("N22").Value = -1840
$ws.Range("H27").Value = 67334390
$ws.Range("I27").Value = 77693340
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 77693340
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = -77693233
$ws.Range("N27").Value = -1464
$ws.Range("H41").Value = 8000
$ws.Range("J41").Value = 8000
$ws.Range("L41").Value = 8000
$ws.Range("N41").Value = -8876
$ws.Range("H61").Value = 13138.5
$ws.Range("I61").Value = 13138.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 13138.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -12936.5
$ws.Range("N61").Value = ""
$ws.Range("H113").Value = 13138.5
$ws.Range("I113").Value = 13138.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 13138.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -10968.5
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 3356.324
$ws.Range("I132").Value = 2867.3447
$ws.Range("J132").Value = 5537.923
$ws.Range("K132").Value = 8602.034100000001
$ws.Range("L132").Value = 16613.769
$ws.Range("M132").Value = -6072.034100000001
$ws.Range("N132").Value = -21673.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H106").Value = 23000
$ws.Range("J106").Value = 23000
$ws.Range("L106").Value = 23000
$ws.Range("N106").Value = -25524
